$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values are written as text (matching the
# workbook's original inline-string cell content) instead of being
# auto-converted to numbers by Excel.
$textCells = @('D5', 'D6', 'D7', 'D8', 'D16', 'D19', 'D20', 'D21', 'D23', 'D24', 'D26', 'D27', 'D29', 'D31', 'D33', 'D36', 'D37', 'D38', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D49', 'D50')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '65.769.71'
$ws.Range('E2').Value = '  +1.92%  '
$ws.Range('D3').Value = '3.484.99'
$ws.Range('E3').Value = '  +0.68%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '580.43'
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').Value = '160.92'
$ws.Range('E6').Value = '  +1.68%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.604'
$ws.Range('E8').Value = '  +8.48%  '
$ws.Range('D9').Value = '3.487.68'
$ws.Range('E9').Value = '  +0.73%  '
$ws.Range('E10').Value = '  -3.56%  '
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').Value = '4.088.31'
$ws.Range('E13').Value = '  +0.65%  '
$ws.Range('E14').Value = '  -1.00%  '
$ws.Range('E15').Value = '  -0.77%  '
$ws.Range('D16').Value = '28.78'
$ws.Range('E16').Value = '  +3.23%  '
$ws.Range('D17').Value = '65.724.36'
$ws.Range('E17').Value = '  +1.78%  '
$ws.Range('D18').Value = '3.515.68'
$ws.Range('E18').Value = '  +1.49%  '
$ws.Range('D19').Value = '6.49'
$ws.Range('E19').Value = '  +0.75%  '
$ws.Range('D20').Value = '14.30'
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').Value = '390.46'
$ws.Range('E21').Value = '  -1.40%  '
$ws.Range('E22').Value = '  -2.97%  '
$ws.Range('D23').Value = '0.552'
$ws.Range('E23').Value = '  +0.97%  '
$ws.Range('D24').Value = '73.59'
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = '0.0000125'
$ws.Range('E26').Value = '  +1.65%  '
$ws.Range('D27').Value = '9.82'
$ws.Range('E27').Value = '  +2.42%  '
$ws.Range('E28').Value = '  -0.50%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('E30').Value = '  +4.82%  '
$ws.Range('D31').Value = '1.44'
$ws.Range('E31').Value = '  +3.86%  '
$ws.Range('E32').Value = '  +1.04%  '
$ws.Range('D33').Value = '23.76'
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('E34').Value = '  -2.78%  '
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('D36').Value = '7.14'
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('D37').Value = '1.55'
$ws.Range('E37').Value = '  +4.97%  '
$ws.Range('D38').Value = '162.85'
$ws.Range('E38').Value = '  +1.40%  '
$ws.Range('E39').Value = '  +4.44%  '
$ws.Range('D40').Value = '3.083.15'
$ws.Range('D41').Value = '0.0774'
$ws.Range('E41').Value = '  -1.52%  '
$ws.Range('D42').Value = '27.23'
$ws.Range('E42').Value = '  -1.49%  '
$ws.Range('D43').Value = '0.0324'
$ws.Range('E43').Value = '  -0.09%  '
$ws.Range('D44').Value = '4.54'
$ws.Range('E44').Value = '  +2.30%  '
$ws.Range('D45').Value = '43.10'
$ws.Range('E45').Value = '  +1.97%  '
$ws.Range('D46').Value = '0.780'
$ws.Range('E46').Value = '  +0.40%  '
$ws.Range('D47').Value = '25.67'
$ws.Range('E47').Value = '  +6.30%  '
$ws.Range('E48').Value = '  +3.06%  '
$ws.Range('D49').Value = '2.25'
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('D50').Value = '6.72'
$ws.Range('E50').Value = '  +2.52%  '
$ws.Range('E51').Value = '  +3.63%  '
